$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 29 ---
$ws.Range("A29").Value = 131237131
$ws.Range("B29").Value = 57881
$ws.Range("D29").Value = "NT"
$ws.Range("E29").Value = 100049
$ws.Range("F29").Value = "Spillkråka"
$ws.Range("G29").Value = "Dryocopus martius"
$ws.Range("H29").Value = "(Linnaeus, 1758)"
$ws.Range("M29").Value = "födosökande"
$ws.Range("P29").Value = "Udden i Råsjön, Ög"
$ws.Range("Q29").Value = 568107
$ws.Range("R29").Value = 6506605
$ws.Range("S29").Value = 10
$ws.Range("T29").Value = "Östergötland"
$ws.Range("U29").Value = "Norrköping"
$ws.Range("V29").Value = "Östergötland"
$ws.Range("W29").Value = "Kvillinge"
$ws.Range("Y29").Value = "'2026-02-19"
$ws.Range("AA29").Value = "'2026-02-19"
$ws.Range("AD29").Value = $false
$ws.Range("AE29").Value = $false
$ws.Range("AG29").Value = $false
$ws.Range("AW29").Value = "Anette Källman"
$ws.Range("AX29").Value = "Anette Källman"

# --- Row 30 ---
$ws.Range("A30").Value = 131237070
$ws.Range("B30").Value = 57881
$ws.Range("D30").Value = "NT"
$ws.Range("E30").Value = 100049
$ws.Range("F30").Value = "Spillkråka"
$ws.Range("G30").Value = "Dryocopus martius"
$ws.Range("H30").Value = "(Linnaeus, 1758)"
$ws.Range("M30").Value = "gammalt bo"
$ws.Range("P30").Value = "Udden i Råsjön, Ög"
$ws.Range("Q30").Value = 568119
$ws.Range("R30").Value = 6506603
$ws.Range("S30").Value = 10
$ws.Range("T30").Value = "Östergötland"
$ws.Range("U30").Value = "Norrköping"
$ws.Range("V30").Value = "Östergötland"
$ws.Range("W30").Value = "Kvillinge"
$ws.Range("Y30").Value = "'2026-02-19"
$ws.Range("AA30").Value = "'2026-02-19"
$ws.Range("AD30").Value = $false
$ws.Range("AE30").Value = $false
$ws.Range("AG30").Value = $false
$ws.Range("AW30").Value = "Anette Källman"
$ws.Range("AX30").Value = "Anette Källman"

Write-Output "rows 29-30 added"
